$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '26.606.71'
$ws.Cells.Item(2, 5).Value = '  -0.15%  '
$ws.Cells.Item(3, 4).Value = '1.644.39'
$ws.Cells.Item(3, 5).Value = '  +0.68%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.01'
$ws.Cells.Item(4, 4).Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.27%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '215.86'
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.31%  '
$ws.Cells.Item(6, 5).Value = '  +0.71%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.00'
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  +0.22%  '
$ws.Cells.Item(8, 5).Value = '  -0.15%  '
$ws.Cells.Item(9, 5).Value = '  +0.77%  '
$ws.Cells.Item(10, 5).Value = '  +0.69%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.0843'
$ws.Cells.Item(11, 4).Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  +0.13%  '
$ws.Cells.Item(12, 4).Value = '1.874.37'
$ws.Cells.Item(12, 5).Value = '  +0.74%  '
$ws.Cells.Item(13, 5).Value = '  +3.26%  '
$ws.Cells.Item(14, 4).Value = '1.645.70'
$ws.Cells.Item(14, 5).Value = '  +1.48%  '
$ws.Cells.Item(15, 5).Value = '  +1.88%  '
$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '66.14'
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +4.50%  '
$ws.Cells.Item(17, 4).Value = '26.664.96'
$ws.Cells.Item(17, 5).Value = '  +0.14%  '
$ws.Cells.Item(18, 5).Value = '  +1.54%  '
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '218.17'
$ws.Cells.Item(19, 4).Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -0.43%  '
$ws.Cells.Item(20, 5).Value = '  +0.31%  '
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '4.37'
$ws.Cells.Item(21, 4).Style = "Normal"
$ws.Cells.Item(21, 5).Value = '  +2.12%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '6.31'
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +1.83%  '
$ws.Cells.Item(23, 5).Value = '  +1.65%  '
$ws.Cells.Item(24, 5).Value = '  +9.70%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '146.61'
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -1.40%  '
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '1.01'
$ws.Cells.Item(26, 4).Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +0.30%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '0.120'
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -0.44%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '7.13'
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  +2.55%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '0.0516'
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  +1.93%  '
$ws.Cells.Item(31, 5).Value = '  +1.16%  '
$ws.Cells.Item(32, 5).Value = '  +3.02%  '
$ws.Cells.Item(33, 5).Value = '  +2.53%  '
$ws.Cells.Item(34, 4).Value = '1.279.44'
$ws.Cells.Item(34, 5).Value = '  +5.71%  '
$ws.Cells.Item(35, 5).Value = '  +2.20%  '
$ws.Cells.Item(36, 5).Value = '  +6.57%  '
$ws.Cells.Item(37, 5).Value = '  +0.17%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.527'
$ws.Cells.Item(38, 4).Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +4.92%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.826'
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  +2.00%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '1.00'
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  +0.22%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '0.807'
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  +1.85%  '
$ws.Cells.Item(42, 5).Value = '  -1.86%  '
$ws.Cells.Item(43, 5).Value = '  +0.72%  '
$ws.Cells.Item(44, 4).Value = '1.785.71'
$ws.Cells.Item(44, 5).Value = '  +0.98%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '93.18'
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +0.25%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '59.68'
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  +9.17%  '
$ws.Cells.Item(47, 5).Value = '  +4.27%  '
$ws.Cells.Item(48, 5).Value = '  +0.81%  '
$ws.Cells.Item(49, 5).Value = '  +2.77%  '
$ws.Cells.Item(50, 5).Value = '  +3.80%  '
$ws.Cells.Item(51, 5).Value = '  -0.60%  '
